$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1056
$ws.Range("I32").Value = 827
$ws.Range("J32").Value = 1285
$ws.Range("K32").Value = 827
$ws.Range("L32").Value = 1285
$ws.Range("M32").Value = -501
$ws.Range("N32").Value = -1937
$ws.Range("H38").Value = 1260.909
$ws.Range("J38").Value = 4432
$ws.Range("L38").Value = 13296
$ws.Range("N38").Value = -14040
$ws.Range("H40").Value = 1799.1111
$ws.Range("I40").Value = 1559
$ws.Range("J40").Value = 1919.1666
$ws.Range("K40").Value = 1559
$ws.Range("L40").Value = 1919.1666
$ws.Range("M40").Value = -1384
$ws.Range("N40").Value = -2269.1666
$ws.Range("H70").Value = 3158.5454
$ws.Range("I70").Value = 990
$ws.Range("J70").Value = 3375.4
$ws.Range("K70").Value = 2970
$ws.Range("L70").Value = 10126.2
$ws.Range("M70").Value = -2700
$ws.Range("N70").Value = -10666.2
$ws.Range("H73").Value = 3158.5454
$ws.Range("I73").Value = 990
$ws.Range("J73").Value = 3375.4
$ws.Range("K73").Value = 2970
$ws.Range("L73").Value = 10126.2
$ws.Range("M73").Value = -2034
$ws.Range("N73").Value = -11998.2
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 3568.3333
$ws.Range("I132").Value = 3765.4546
$ws.Range("K132").Value = 11296.3638
$ws.Range("M132").Value = -8766.363799999999
$ws.Range("H135").Value = 946.75
$ws.Range("I135").Value = 948.2857
$ws.Range("K135").Value = 8534.5713
$ws.Range("M135").Value = -5999.5713
$ws.Range("H138").Value = 4244.4
$ws.Range("J138").Value = 4244.4
$ws.Range("L138").Value = 12733.2
$ws.Range("N138").Value = -23013.2
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 2320.4375
$ws.Range("I141").Value = 2320.4375
$ws.Range("K141").Value = 6961.3125
$ws.Range("M141").Value = -1781.3125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1959.64
$ws.Range("I32").Value = 1959.64
$ws.Range("K32").Value = 1959.64
$ws.Range("M32").Value = -1672.64

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1866.5
$ws.Range("I5").Value = 1299.75
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 1299.75
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1186.75
$ws.Range("N5").Value = -3226
$ws.Range("H95").Value = 12624
$ws.Range("J95").Value = 12624
$ws.Range("L95").Value = 12624
$ws.Range("N95").Value = -18116

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 367.2857
$ws.Range("I2").Value = 386.83334
$ws.Range("K2").Value = 386.83334
$ws.Range("M2").Value = -273.83334
$ws.Range("H28").Value = 11183
$ws.Range("J28").Value = 11183
$ws.Range("L28").Value = 11183
$ws.Range("N28").Value = -11673
$ws.Range("H43").Value = 10499.8
$ws.Range("J43").Value = 10499.8
$ws.Range("L43").Value = 10499.8
$ws.Range("N43").Value = -10867.8
$ws.Range("H54").Value = 4538.8887
$ws.Range("I54").Value = 2000
$ws.Range("J54").Value = 9616.666999999999
$ws.Range("K54").Value = 2000
$ws.Range("L54").Value = 9616.666999999999
$ws.Range("M54").Value = -1342
$ws.Range("N54").Value = -10932.667
$ws.Range("H58").Value = 5608.1665
$ws.Range("I58").Value = 5218
$ws.Range("J58").Value = 9900
$ws.Range("K58").Value = 5218
$ws.Range("L58").Value = 9900
$ws.Range("M58").Value = -5015
$ws.Range("N58").Value = -10306
$ws.Range("H86").Value = 9424.75
$ws.Range("I86").Value = 8689.888999999999
$ws.Range("K86").Value = 8689.888999999999
$ws.Range("M86").Value = -7566.888999999999
$ws.Range("H89").Value = 9424.75
$ws.Range("I89").Value = 8689.888999999999
$ws.Range("K89").Value = 43449.44499999999
$ws.Range("M89").Value = -37833.44499999999
$ws.Range("H101").Value = 10499.8
$ws.Range("J101").Value = 10499.8
$ws.Range("L101").Value = 10499.8
$ws.Range("N101").Value = -16989.8
$ws.Range("H134").Value = 57680.445
$ws.Range("I134").Value = 72446.42999999999
$ws.Range("J134").Value = 5999.5
$ws.Range("K134").Value = 217339.29
$ws.Range("L134").Value = 17998.5
$ws.Range("M134").Value = -214804.29
$ws.Range("N134").Value = -23068.5
$ws.Range("H136").Value = 5608.1665
$ws.Range("I136").Value = 5218
$ws.Range("J136").Value = 9900
$ws.Range("K136").Value = 15654
$ws.Range("L136").Value = 29700
$ws.Range("M136").Value = -13104
$ws.Range("N136").Value = -34800

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("K3").Value = 15000
$ws.Range("M3").Value = -14888
$ws.Range("H7").Value = 312.14285
$ws.Range("I7").Value = 250.6
$ws.Range("J7").Value = 346.33334
$ws.Range("K7").Value = 751.8
$ws.Range("L7").Value = 1039.00002
$ws.Range("M7").Value = -639.8
$ws.Range("N7").Value = -1263.00002
$ws.Range("H17").Value = 1065.6666
$ws.Range("I17").Value = 398.75
$ws.Range("J17").Value = 1599.2
$ws.Range("K17").Value = 1196.25
$ws.Range("L17").Value = 4797.6
$ws.Range("M17").Value = -1027.25
$ws.Range("N17").Value = -5135.6
$ws.Range("H36").Value = 316.66666
$ws.Range("I36").Value = 325
$ws.Range("K36").Value = 975
$ws.Range("M36").Value = -806
$ws.Range("H39").Value = 5000
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H128").Value = 339751.44
$ws.Range("I128").Value = 339751.44
$ws.Range("K128").Value = 1019254.32
$ws.Range("M128").Value = -1014274.32
$ws.Range("H129").Value = 1400.8334
$ws.Range("J129").Value = 5997
$ws.Range("L129").Value = 17991
$ws.Range("N129").Value = -27991
$ws.Range("H131").Value = 1195
$ws.Range("J131").Value = 1500
$ws.Range("L131").Value = 4500
$ws.Range("N131").Value = -14580

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3873.25
$ws.Range("I70").Value = 3499
$ws.Range("K70").Value = 3499
$ws.Range("M70").Value = -3229
$ws.Range("H73").Value = 3873.25
$ws.Range("I73").Value = 3499
$ws.Range("K73").Value = 3499
$ws.Range("M73").Value = -2563
$ws.Range("H80").Value = 3367
$ws.Range("I80").Value = 2741.6667
$ws.Range("J80").Value = 3992.3333
$ws.Range("K80").Value = 2741.6667
$ws.Range("L80").Value = 3992.3333
$ws.Range("M80").Value = -1743.6667
$ws.Range("N80").Value = -5988.3333
$ws.Range("H83").Value = 3367
$ws.Range("I83").Value = 2741.6667
$ws.Range("J83").Value = 3992.3333
$ws.Range("K83").Value = 13708.3335
$ws.Range("L83").Value = 19961.6665
$ws.Range("M83").Value = -8716.333500000001
$ws.Range("N83").Value = -29945.6665

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2250
$ws.Range("J7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("N7").Value = -2724
$ws.Range("H126").Value = 2250
$ws.Range("J126").Value = 2500
$ws.Range("L126").Value = 7500
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 6884.0625
$ws.Range("I132").Value = 6343.857
$ws.Range("J132").Value = 7304.222
$ws.Range("K132").Value = 19031.571
$ws.Range("L132").Value = 21912.666
$ws.Range("M132").Value = -16501.571
$ws.Range("N132").Value = -26972.666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 48000
$ws.Range("J112").Value = 48000
$ws.Range("L112").Value = 48000
$ws.Range("N112").Value = -50954
$ws.Range("H113").Value = 803.4
$ws.Range("I113").Value = 775.25
$ws.Range("K113").Value = 2325.75
$ws.Range("M113").Value = -155.75
